$d = $word.ActiveDocument

$newText = "Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi himmeimpiä näkyvissä olevia tähtiä keinona mitata valonsaastetta tietyssä paikassa. Paikallistamalla ja tarkkailemalla Leijonan tähtikuvio miten valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. Antamasi tiedot päivittyvät heti verkossa olevaan tietokantaan, ja näin saadaan käsitys siitä minkä verran taivaan tähdistä on missäkin nähtävissä."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Osallistut maailmanlaajuiseen tapahtumaan*") {
        $rng = $p.Range
        $rng.MoveEnd(1, -1)
        $rng.Delete()

        $rng2 = $p.Range
        $rng2.MoveEnd(1, -1)
        $rng2.InsertAfter($newText)
        break
    }
}
